$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text, if assigned directly, would be auto-converted by Excel
# into a number (e.g. "1.0000" -> 1, "110.00" -> 110), losing the exact display
# text captured in the source data. For these we temporarily force the cell to
# Text format, assign the literal string, then restore the cell to the default
# "Normal" style so no stray formatting is left behind on the cell.
$textForcedUpdates = @{
    "D4" = '0.9997'
    "D5" = '329.22'
    "D6" = '1.0000'
    "D7" = '0.4564'
    "D8" = '0.3938'
    "D9" = '47.56'
    "D10" = '0.07829'
    "D11" = '0.9827'
    "D12" = '21.32'
    "D14" = '5.828'
    "D15" = '6.963'
    "D17" = '87.99'
    "D18" = '0.06514'
    "D20" = '17.10'
    "D23" = '5.300'
    "D24" = '10.81'
    "D25" = '2.253'
    "D27" = '157.01'
    "D28" = '19.21'
    "D29" = '2.061'
    "D30" = '5.297'
    "D31" = '116.38'
    "D32" = '0.9416'
    "D33" = '0.09276'
    "D34" = '3.602'
    "D35" = '1.384'
    "D36" = '5.207'
    "D37" = '0.06006'
    "D38" = '0.02200'
    "D39" = '8.248'
    "D40" = '1.160'
    "D41" = '1.000'
    "D42" = '0.5693'
    "D43" = '0.1797'
    "D44" = '9.968'
    "D45" = '1.257'
    "D46" = '2.309'
    "D47" = '0.5413'
    "D48" = '11.82'
    "D49" = '0.07160'
    "D50" = '1.868'
    "D51" = '110.00'
}

# Cells whose new text is not a valid number literal (e.g. thousand-grouped
# "28.318.65" or percentages with surrounding spaces) and so Excel keeps them
# as plain text automatically; these can be assigned directly.
$plainUpdates = @{
    "D2" = '28.318.65'
    "E2" = '  +0.27%  '
    "D3" = '1.858.64'
    "E3" = '  -0.61%  '
    "E4" = '  -0.05%  '
    "E5" = '  -2.33%  '
    "E6" = '  -0.05%  '
    "E7" = '  -2.88%  '
    "E8" = '  +0.41%  '
    "E9" = '  +1.03%  '
    "E10" = '  -1.78%  '
    "E11" = '  -2.74%  '
    "E12" = '  -1.89%  '
    "D13" = '1.847.60'
    "E13" = '  -1.27%  '
    "E14" = '  -2.80%  '
    "E15" = '  -4.19%  '
    "E16" = '  +0.04%  '
    "E17" = '  -3.64%  '
    "E18" = '  -1.19%  '
    "E19" = '  -2.44%  '
    "E20" = '  -3.56%  '
    "E21" = '  +0.03%  '
    "D22" = '28.281.52'
    "E22" = '  +0.14%  '
    "E23" = '  -2.76%  '
    "E24" = '  -2.37%  '
    "E25" = '  -1.80%  '
    "D26" = '2.072.24'
    "E26" = '  -0.95%  '
    "E27" = '  -1.40%  '
    "E28" = '  -4.11%  '
    "E29" = '  -3.96%  '
    "E30" = '  -3.70%  '
    "E31" = '  -2.87%  '
    "E32" = '  -3.72%  '
    "E33" = '  -1.89%  '
    "E34" = '  +0.74%  '
    "E35" = '  +0.68%  '
    "E36" = '  -2.74%  '
    "E37" = '  -1.61%  '
    "E38" = '  -2.84%  '
    "E39" = '  -2.58%  '
    "E40" = '  -1.64%  '
    "E41" = '  +0.05%  '
    "E42" = '  -4.84%  '
    "E43" = '  -4.61%  '
    "E44" = '  -4.16%  '
    "E45" = '  -2.90%  '
    "E46" = '  +17.05%  '
    "E47" = '  -3.75%  '
    "E48" = '  -3.39%  '
    "E49" = '  +4.07%  '
    "E50" = '  -5.46%  '
    "E51" = '  -1.17%  '
}

foreach ($ref in $textForcedUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$ref]
    $cell.Style = "Normal"
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}
